$d = $word.ActiveDocument

# --- Change 1: "Person" -> "c" in the table cell ---
$d.Content.Find.Execute("Person", $false, $false, $false, $false, $false, $true, 1, $false, "c", 2)

# --- Change 2: append four new bold paragraphs after the trailing "}" ---
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last.Range
$p1.Text = "Answer :- "
$p1.Font.Bold = $true
$p1.Font.BoldBi = $true
$p1.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last.Range
$p2.Text = "true"
$p2.Font.Bold = $true
$p2.Font.BoldBi = $true
$p2.InsertParagraphAfter()

$p3 = $d.Paragraphs.Last.Range
$p3.Text = "true"
$p3.Font.Bold = $true
$p3.Font.BoldBi = $true
$p3.InsertParagraphAfter()

$p4 = $d.Paragraphs.Last.Range
$p4.Text = "true"
$p4.Font.Bold = $true
$p4.Font.BoldBi = $true

Write-Output "edits applied"
